# Update cryptos list (price + 1h volume change) with refreshed figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.833.31'
$ws.Range("E2").Value = '  +3.49%  '
$ws.Range("D3").Value = '3.265.67'
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'580.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("D6").Value = "'181.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.16%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '3.263.91'
$ws.Range("E9").Value = '  +2.98%  '
$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.58%  '
$ws.Range("E11").Value = '  +2.80%  '
$ws.Range("E12").Value = '  +6.27%  '
$ws.Range("D13").Value = '3.824.65'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").Value = "'28.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.25%  '
$ws.Range("D16").Value = '67.747.45'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("D17").Value = "'0.0000169"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = '3.245.93'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("D19").Value = "'5.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").Value = "'13.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.07%  '
$ws.Range("D21").Value = "'375.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.87%  '
$ws.Range("D22").Value = "'7.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.17%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = "'71.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").Value = "'0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.70%  '
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.58%  '
$ws.Range("D27").Value = "'9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +3.04%  '
$ws.Range("D31").Value = "'5.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.39%  '
$ws.Range("D32").Value = "'22.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.62%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +5.60%  '
$ws.Range("E35").Value = '  +4.83%  '
$ws.Range("E36").Value = '  +4.60%  '
$ws.Range("D37").Value = "'163.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("D38").Value = "'0.851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.24%  '
$ws.Range("D39").Value = "'1.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("D40").Value = "'6.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.16%  '
$ws.Range("D41").Value = "'26.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("E42").Value = '  +10.78%  '
$ws.Range("D43").Value = "'2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.06%  '
$ws.Range("D44").Value = '2.694.63'
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("D45").Value = "'351.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.35%  '
$ws.Range("D46").Value = "'25.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.17%  '
$ws.Range("D47").Value = "'40.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("D49").Value = "'0.0281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.55%  '
$ws.Range("E51").Value = '  +0.24%  '
